$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2:E3").Value = "2016-03-19 04:10:19"
$wsZh.Range("H2:H3").Value = "2016-03-19 04:11:01"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2:E3").Value = "2016-03-19 04:10:27"
$wsDe.Range("H2:H3").Value = "2016-03-19 04:11:15"
